# Generate Report for Handback
# Updates the localization-status workbook to reflect a completed handback:
#  - Status column moves from "Ready for handoff" to "Handed back: in sync with en-US"
#    (this text is shared across the Overview summary columns and the per-language
#    sheets, so every cell showing the old status gets updated)
#  - zh-cn / de-de sheets get their "Latest Target File" / "Latest Handback File" /
#    "Latest Handback DateTime" columns populated
#  - Column widths widen to fit the new, longer status text
#  - A hyperlink (matching the source-file hyperlink in column A) is added on the
#    "Latest Target File" cell for each data row

$wb = $excel.ActiveWorkbook

$baseUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/605c39174179330d95862132aa9635ca6ac3226f/e2e/"

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"
$zhDateTime = "2016-09-06 05:19:43"
$deDateTime = "2016-09-06 05:19:51"

# ---- Overview sheet: widen the zh-cn / de-de status columns and refresh status text ----
$ov = $wb.Worksheets.Item("Overview")
$ov.Columns.Item(5).ColumnWidth = 29.9777047293527
$ov.Columns.Item(6).ColumnWidth = 29.9777047293527

for ($row = 2; $row -le 3; $row++) {
    $ov.Cells.Item($row, 5).Value = $newStatus
    $ov.Cells.Item($row, 6).Value = $newStatus
}

# ---- Per-language sheets ----
$langs = @(
    @{ Sheet = "zh-cn"; File = "462b48b3-567d-4863-b338-db92d04a245e.afeafe056b3d5c2feb2964de8ab7465eacddbd50.zh-cn.xlf"; DateTime = $zhDateTime },
    @{ Sheet = "de-de"; File = "462b48b3-567d-4863-b338-db92d04a245e.afeafe056b3d5c2feb2964de8ab7465eacddbd50.de-de.xlf"; DateTime = $deDateTime }
)

foreach ($lang in $langs) {
    $ws = $wb.Worksheets.Item($lang.Sheet)

    # Widen Status (C), Latest Target File (I) and Latest Handback File (J) columns
    $ws.Columns.Item(3).ColumnWidth = 29.9777047293527
    $ws.Columns.Item(9).ColumnWidth = 40
    $ws.Columns.Item(10).ColumnWidth = 40

    # Both data rows point their "Latest Target File" hyperlink at the same
    # 462b48b3-... source document (matches the handback report's own output).
    $targetSourceName = "462b48b3-567d-4863-b338-db92d04a245e.md"
    $target = $baseUrl + $targetSourceName

    for ($row = 2; $row -le 3; $row++) {
        # Status column -> handed back
        $ws.Cells.Item($row, 3).Value = $newStatus

        # Latest Target File (I) becomes a hyperlink to the source doc
        $ws.Hyperlinks.Add($ws.Cells.Item($row, 9), $target, "", "", $targetSourceName)

        # Latest Handback File (J) = the xlf file that was handed back
        $ws.Cells.Item($row, 10).Value = $lang.File

        # Latest Handback DateTime (K)
        $ws.Cells.Item($row, 11).Value = $lang.DateTime
    }
}
